# Edit script: rewrite the closing "Hello there..." paragraph into four
# paragraphs (an intro blurb, a styled pull-quote, and two more narrative
# paragraphs), per the target diff.

$d = $word.ActiveDocument

# Locate the target paragraph robustly by its distinctive text, rather than
# assuming a fixed index.
$targetIdx = -1
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    $t = $para.Range.Text
    if ($t -like 'Hello there*' -and $t -like '*amateur scientist*') {
        $targetIdx = $i
    }
}

if ($targetIdx -eq -1) {
    throw "Could not locate the target 'Hello there...' paragraph"
}

$target = $d.Paragraphs.Item($targetIdx)

# Replace the whole paragraph's text (all of its runs) with the new intro
# text in one shot, using a Range built from explicit offsets -- this
# collapses the old multi-run paragraph down to a single run.
$r1 = $d.Range($target.Range.Start, $target.Range.End)
$r1.Text = 'Hello there! I’m Brandon Hawley a dedicated father of four, expert software architect, aspiring scientist, artist, and game enthusiast. I believe creativity and imagination are some of our greatest gifts as humans; my passion in life derives from this very idea. I dream to further human development, improve the human condition, write narratives and create works of art that inspire humanity.'

# Insert three more (currently empty) paragraphs right after it -- do this
# before applying any custom formatting so the new paragraphs don't inherit
# direct formatting from a differently-styled neighbor. Inserting after the
# same anchor range four times in a row stacks the new empty paragraphs
# immediately following it (and leaves one extra blank paragraph behind,
# which gets deleted below once we know where it landed).
$r1.InsertParagraphAfter()
$r1.InsertParagraphAfter()
$r1.InsertParagraphAfter()
$r1.InsertParagraphAfter()

$extra = $d.Paragraphs.Item($targetIdx + 4)
$extra.Range.Delete()

# Paragraph 1 (the intro blurb): indent the first line half an inch.
$p1Para = $d.Paragraphs.Item($targetIdx)
$p1Para.Range.ParagraphFormat.FirstLineIndent = 36

# Paragraph 2: the styled pull-quote -- bold, italic, grey (accent3-ish),
# indented from the left.
$p2Para = $d.Paragraphs.Item($targetIdx + 1)
$p2Para.Range.Text = '“Intelligence without creativity is lost.”'
$p2Para.Range.ParagraphFormat.LeftIndent = 108
$p2Para.Range.Font.Bold = 1
$p2Para.Range.Font.Italic = 1
$p2Para.Range.Font.ItalicBi = 1
$p2Para.Range.Font.Color = 10855845

# Paragraph 3: continues the narrative, first-line indented.
$p3Para = $d.Paragraphs.Item($targetIdx + 2)
$p3Para.Range.Text = 'Equipped with a rich tapestry of skills, I''ve ventured across many diverse landscapes of software development. Vanquishing C++, C#, Java, Python, and more, my creative pulse thrives with challenge. The greatest tool I have But my real love comes from creating crafting intricate narratives and immersive experiences using Unreal Engine and Unity Engine. The symphony of challenges posed by these languages fuels my passion, driving me to conquer every intricacy they offer. As I wield these tools, my palette expands to include Blender for captivating visual creations. While my journey includes web development prowess in HTML, CSS, and JavaScript, it''s my love affair with Unreal Engine that takes center stage, breathing life into my now tangible dreams. Igniting my fervor for crafting interactive worlds.'
$p3Para.Range.ParagraphFormat.FirstLineIndent = 36

# Paragraph 4: closes out the narrative, first-line indented.
$p4Para = $d.Paragraphs.Item($targetIdx + 3)
$p4Para.Range.Text = 'With a keyboard as my trusty sword, I embark on magnificent adventures in the worlds of code. Each castle, the challenges I face for every new language. The dragons above, the errors that seem mythically elusive. But, as my skill tree grows, the expansiveness of new worlds open their doors. What used to challenge become mere tools in my arsenal; the languages, like languages of ancients, empower me to craft experiences that captivate and amaze.'
$p4Para.Range.ParagraphFormat.FirstLineIndent = 36

Write-Output "Rewrote paragraph $targetIdx into 4 paragraphs; doc now has $($d.Paragraphs.Count) paragraphs."
